$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (corine-land-cover-2000-nivel-2-descripcion) switches from a
# curated "dimension" to a curated "measure".
$ws.Range("D2").Value = "iaest-measure:corine-land-cover-2000-nivel-2-descripcion"
$ws.Range("D3").Value = "medida"
$ws.Range("D4").Value = "xsd:int"
# Measures do not carry a mapping file, so the old mapping reference is removed.
$ws.Range("D5").Clear()

# Column F (municipio-nombre) switches from a curated "measure" to a
# curated "dimension" referencing an area, with its own URI type.
$ws.Range("F2").Value = "sdmx-dimension:refArea"
$ws.Range("F3").Value = "dim"
$ws.Range("F4").Value = "URI-Municipio"
